# room_schedule_Computer_B-203A.xlsx edit script
# Adds a Monday / Thursday course grid (room "B-203A") to the existing
# time-slot schedule sheet.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Colors (Excel BGR-packed integers, equivalent to RGB(r,g,b))
#   255       -> FF0000 red    (day-of-week header band)
#   65535     -> FFFF00 yellow (room-name band)
#   16755575  -> 77ABFF blue   (course entry band)

# =======================================================================
# Row 1 / Row 2 header blocks: merge the cell pairs first (so the shared
# border isn't split by the merge), then format every cell belonging to
# the same visual style together so the engine can fold the progressive
# style changes into a single shared cell format.
# =======================================================================
$ws.Range("B1:C1").Merge()
$ws.Range("D1:E1").Merge()
$ws.Range("B2:C2").Merge()
$ws.Range("D2:E2").Merge()

$b1 = $ws.Range("B1")
$c1 = $ws.Range("C1")
$d1 = $ws.Range("D1")
$e1 = $ws.Range("E1")
$b2 = $ws.Range("B2")
$c2 = $ws.Range("C2")
$d2 = $ws.Range("D2")
$e2 = $ws.Range("E2")

$b1.Value = "Monday"
$d1.Value = "Thursday"
$b2.Value = "B-203A"
$d2.Value = "B-203A"

# ---- Row 1: bold, red fill, bordered, centered ----
$b1.Font.Bold = $true
$b1.Interior.Color = 255
$b1.Borders.LineStyle = 1
$b1.HorizontalAlignment = -4108
$b1.VerticalAlignment = -4108

$c1.Font.Bold = $true
$c1.Interior.Color = 255
$c1.Borders.LineStyle = 1
$c1.HorizontalAlignment = -4108
$c1.VerticalAlignment = -4108

$d1.Font.Bold = $true
$d1.Interior.Color = 255
$d1.Borders.LineStyle = 1
$d1.HorizontalAlignment = -4108
$d1.VerticalAlignment = -4108

$e1.Font.Bold = $true
$e1.Interior.Color = 255
$e1.Borders.LineStyle = 1
$e1.HorizontalAlignment = -4108
$e1.VerticalAlignment = -4108

# ---- Row 2: bold, yellow fill, bordered, centered ----
$b2.Font.Bold = $true
$b2.Interior.Color = 65535
$b2.Borders.LineStyle = 1
$b2.HorizontalAlignment = -4108
$b2.VerticalAlignment = -4108

$c2.Font.Bold = $true
$c2.Interior.Color = 65535
$c2.Borders.LineStyle = 1
$c2.HorizontalAlignment = -4108
$c2.VerticalAlignment = -4108

$d2.Font.Bold = $true
$d2.Interior.Color = 65535
$d2.Borders.LineStyle = 1
$d2.HorizontalAlignment = -4108
$d2.VerticalAlignment = -4108

$e2.Font.Bold = $true
$e2.Interior.Color = 65535
$e2.Borders.LineStyle = 1
$e2.HorizontalAlignment = -4108
$e2.VerticalAlignment = -4108

# =======================================================================
# Course entries: plain (non-bold) cells with a light-blue fill.
# Monday block (9:00-10:00, rows 7-8) and Thursday block
# (16:30-17:30, rows 31-32).
# =======================================================================
$b7 = $ws.Range("B7")
$c7 = $ws.Range("C7")
$b8 = $ws.Range("B8")
$c8 = $ws.Range("C8")
$d31 = $ws.Range("D31")
$e31 = $ws.Range("E31")
$d32 = $ws.Range("D32")
$e32 = $ws.Range("E32")

$b7.Value = "A4"
$c7.Value = "WTL"
$b8.Value = "TEA"
$c8.Value = "YVD"
$d31.Value = "C2"
$e31.Value = "MIL"
$d32.Value = "SEC"
$e32.Value = "VAMI"

$b7.Interior.Color = 16755575
$b7.Borders.LineStyle = 1
$b7.HorizontalAlignment = -4108
$b7.VerticalAlignment = -4108

$c7.Interior.Color = 16755575
$c7.Borders.LineStyle = 1
$c7.HorizontalAlignment = -4108
$c7.VerticalAlignment = -4108

$b8.Interior.Color = 16755575
$b8.Borders.LineStyle = 1
$b8.HorizontalAlignment = -4108
$b8.VerticalAlignment = -4108

$c8.Interior.Color = 16755575
$c8.Borders.LineStyle = 1
$c8.HorizontalAlignment = -4108
$c8.VerticalAlignment = -4108

$d31.Interior.Color = 16755575
$d31.Borders.LineStyle = 1
$d31.HorizontalAlignment = -4108
$d31.VerticalAlignment = -4108

$e31.Interior.Color = 16755575
$e31.Borders.LineStyle = 1
$e31.HorizontalAlignment = -4108
$e31.VerticalAlignment = -4108

$d32.Interior.Color = 16755575
$d32.Borders.LineStyle = 1
$d32.HorizontalAlignment = -4108
$d32.VerticalAlignment = -4108

$e32.Interior.Color = 16755575
$e32.Borders.LineStyle = 1
$e32.HorizontalAlignment = -4108
$e32.VerticalAlignment = -4108
